# [AFG] added final excel sheets for Afghanistan
# Adds a new worksheet "ODI Batting Extra" after the existing sheets,
# containing per-match batting extras for the player.

$wb = $excel.ActiveWorkbook

# Add the new worksheet at the end of the workbook (after the last existing sheet)
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "ODI Batting Extra"

# Helper to write a value as text (not auto-converted to a number),
# while keeping the cell's style back to the sheet's default ("Normal").
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Clone the bold/bordered header formatting used on the other sheets,
# then overwrite with this sheet's header captions.
$wb.Worksheets.Item("ODI Batting").Range("A1:F1").Copy($ws.Range("A1:F1"))

$ws.Range("A1").Value = "MATCH_CODE"
$ws.Range("B1").Value = "BATTING_POSITION"
$ws.Range("C1").Value = "NUM_4"
$ws.Range("D1").Value = "NUM_6"
$ws.Range("E1").Value = "PERCENT_RUNS_OF_TOTAL"
$ws.Range("F1").Value = "MAN_OF_MATCH"

# Row 2 - match 4530
Set-TextValue $ws.Range("A2") "4530"
$ws.Range("B2").Value = 2
Set-TextValue $ws.Range("C2") "3"
Set-TextValue $ws.Range("D2") "1"
Set-TextValue $ws.Range("E2") "19.69%"
$ws.Range("F2").Value = "NO"

# Row 3 - match 4538 (only MATCH_CODE and MAN_OF_MATCH known)
Set-TextValue $ws.Range("A3") "4538"
$ws.Range("F3").Value = "NO"

# Row 4 - match 4539
Set-TextValue $ws.Range("A4") "4539"
$ws.Range("B4").Value = 2
Set-TextValue $ws.Range("C4") "4"
Set-TextValue $ws.Range("D4") "1"
Set-TextValue $ws.Range("E4") "18.13%"
$ws.Range("F4").Value = "NO"
